$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.417.85"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "3.490.32"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.55"
$ws.Range("E5").Value = "  +4.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.72"
$ws.Range("E6").Value = "  +8.54%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.490.75"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "4.092.21"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.88"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "66.312.31"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").Value = "3.476.36"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.75"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.00"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.06"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  +5.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +4.58%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.34"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.903"
$ws.Range("E38").Value = "  +9.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0748"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.27"
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.42"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.68"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.59"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.811.33"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.33"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0313"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.51"
$ws.Range("E48").Value = "  +6.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "349.84"
$ws.Range("E49").Value = "  +8.14%  "
$ws.Range("E50").Value = "  +6.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.83"
$ws.Range("E51").Value = "  +8.54%  "
